$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking strings are not
# auto-converted to numbers (matches the inlineStr type in the source file).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.023.33"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "3.083.32"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "578.19"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").Value = "168.57"
$ws.Range("E6").Value = "  -2.77%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.080.36"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").Value = "6.42"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").Value = "0.150"
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").Value = "0.472"
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").Value = "0.0000241"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D14").Value = "36.19"
$ws.Range("E14").Value = "  -2.41%  "
$ws.Range("E15").Value = "  -2.09%  "
$ws.Range("D16").Value = "3.596.64"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").Value = "66.970.07"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "7.01"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").Value = "16.60"
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("D20").Value = "3.080.42"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "490.03"
$ws.Range("E21").Value = "  +3.01%  "
$ws.Range("D22").Value = "7.75"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "0.687"
$ws.Range("E23").Value = "  -3.10%  "
$ws.Range("D24").Value = "82.80"
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("D25").Value = "12.87"
$ws.Range("E25").Value = "  -3.42%  "
$ws.Range("D26").Value = "2.23"
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("D27").Value = "10.25"
$ws.Range("E27").Value = "  +3.79%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("D30").Value = "2.31"
$ws.Range("E30").Value = "  -4.24%  "
$ws.Range("D31").Value = "2.63"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").Value = "27.88"
$ws.Range("E32").Value = "  -2.51%  "
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("D34").Value = "0.0₃0911"
$ws.Range("E34").Value = "  -3.57%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "5.70"
$ws.Range("E36").Value = "  -2.32%  "
$ws.Range("E37").Value = "  -2.23%  "
$ws.Range("D38").Value = "46.58"
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("D39").Value = "0.123"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("D40").Value = "1.99"
$ws.Range("E40").Value = "  -3.93%  "
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("D42").Value = "8.33"
$ws.Range("E42").Value = "  -2.80%  "
$ws.Range("D43").Value = "2.778.82"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").Value = "371.35"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("D45").Value = "0.0345"
$ws.Range("E45").Value = "  -2.32%  "
$ws.Range("D46").Value = "135.64"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").Value = "2.48"
$ws.Range("E47").Value = "  -2.35%  "
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").Value = "24.54"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("E51").Value = "  -1.19%  "
